$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to be inserted into column B for each row (2..16),
# pushing the existing B:K values one column to the right (and
# dropping whatever value falls off the end of column K).
$newValues = @{
    2  = -0.6603092772102132
    3  = -0.15162438770796
    4  = -0.2053460154962278
    5  = 0.6162032393936197
    6  = 1.652643173475852
    7  = 0.3110387314724781
    8  = 0.2388379152847414
    9  = 0.6508000635779043
    10 = 0.2387740594105157
    11 = 0.3465902496671606
    12 = 0.00230005330798793
    13 = -0.1902738424076751
    14 = -0.3325070745318338
    15 = 0.1656141382254278
    16 = -0.09587373626955231
}

foreach ($row in 2..16) {
    # Shift existing values in columns B..K one column to the right,
    # working from the rightmost column down to column C so that
    # values are not clobbered before they are copied. The old
    # column K value (if any) simply falls off the end.
    for ($col = 11; $col -ge 3; $col--) {
        $srcCell = $ws.Cells.Item($row, $col - 1)
        $dstCell = $ws.Cells.Item($row, $col)
        $dstCell.Value = $srcCell.Value()
    }

    # Write the new value into column B for this row.
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}
